$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 984.43475
$ws.Range("I98").Value = 1040
$ws.Range("J98").Value = 614
$ws.Range("K98").Value = 1040
$ws.Range("L98").Value = 614
$ws.Range("M98").Value = 458
$ws.Range("N98").Value = -3610

$ws.Range("H100").Value = 2456.6667
$ws.Range("I100").Value = 2585
$ws.Range("J100").Value = 2200
$ws.Range("K100").Value = 2585
$ws.Range("L100").Value = 2200
$ws.Range("M100").Value = -2044
$ws.Range("N100").Value = -3282

$ws.Range("H116").Value = 2755.5557
$ws.Range("I116").Value = 9000
$ws.Range("J116").Value = 2388.2354
$ws.Range("K116").Value = 9000
$ws.Range("L116").Value = 2388.2354
$ws.Range("M116").Value = -5558
$ws.Range("N116").Value = -9272.2354

$ws.Range("H122").Value = 984.43475
$ws.Range("I122").Value = 1040
$ws.Range("J122").Value = 614
$ws.Range("K122").Value = 3120
$ws.Range("L122").Value = 1842
$ws.Range("M122").Value = -670
$ws.Range("N122").Value = -6742

$ws.Range("H125").Value = 1987
$ws.Range("I125").Value = 1861
$ws.Range("J125").Value = 2081.5
$ws.Range("K125").Value = 16749
$ws.Range("L125").Value = 18733.5
$ws.Range("M125").Value = -14289
$ws.Range("N125").Value = -23653.5

$ws.Range("H138").Value = 1873.79
$ws.Range("J138").Value = 2184.6711
$ws.Range("L138").Value = 6554.013300000001
$ws.Range("N138").Value = -16834.0133

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11124968
$ws.Range("I32").Value = 15391895
$ws.Range("K32").Value = 15391895
$ws.Range("M32").Value = -15391608

$ws.Range("H74").Value = 1472.3928
$ws.Range("I74").Value = 1108.52
$ws.Range("J74").Value = 4504.6665
$ws.Range("K74").Value = 1108.52
$ws.Range("L74").Value = 4504.6665
$ws.Range("M74").Value = -234.52
$ws.Range("N74").Value = -6252.6665

$ws.Range("H77").Value = 1472.3928
$ws.Range("I77").Value = 1108.52
$ws.Range("J77").Value = 4504.6665
$ws.Range("K77").Value = 5542.6
$ws.Range("L77").Value = 22523.3325
$ws.Range("M77").Value = -1174.6
$ws.Range("N77").Value = -31259.3325

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 5450.4165
$ws.Range("I7").Value = 2500
$ws.Range("J7").Value = 5718.636
$ws.Range("K7").Value = 2500
$ws.Range("L7").Value = 5718.636
$ws.Range("M7").Value = -2387
$ws.Range("N7").Value = -5944.636

$ws.Range("H99").Value = 2490.182
$ws.Range("I99").Value = 1745
$ws.Range("K99").Value = 1745
$ws.Range("M99").Value = -247

$ws.Range("H105").Value = 2015.3158
$ws.Range("I105").Value = 1635.2307
$ws.Range("J105").Value = 2838.8333
$ws.Range("K105").Value = 1635.2307
$ws.Range("L105").Value = 2838.8333
$ws.Range("M105").Value = 111.7692999999999
$ws.Range("N105").Value = -6332.8333

$ws.Range("H134").Value = 2277.8857
$ws.Range("I134").Value = 1735.3334
$ws.Range("J134").Value = 3091.7144
$ws.Range("K134").Value = 5206.0002
$ws.Range("L134").Value = 9275.143199999999
$ws.Range("M134").Value = -2671.0002
$ws.Range("N134").Value = -14345.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1789.4736
$ws.Range("I99").Value = 1638.4615
$ws.Range("J99").Value = 2116.6667
$ws.Range("K99").Value = 1638.4615
$ws.Range("L99").Value = 2116.6667
$ws.Range("M99").Value = -140.4614999999999
$ws.Range("N99").Value = -5112.6667

$ws.Range("H105").Value = 822
$ws.Range("I105").Value = 703.3333
$ws.Range("J105").Value = 1000
$ws.Range("K105").Value = 703.3333
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 1043.6667
$ws.Range("N105").Value = -4494

$ws.Range("H126").Value = 1789.4736
$ws.Range("I126").Value = 1638.4615
$ws.Range("J126").Value = 2116.6667
$ws.Range("K126").Value = 4915.3845
$ws.Range("L126").Value = 6350.000100000001
$ws.Range("M126").Value = -2445.3845
$ws.Range("N126").Value = -11290.0001

$ws.Range("H132").Value = 1876.3334
$ws.Range("I132").Value = 1461.1666
$ws.Range("J132").Value = 2374.5334
$ws.Range("K132").Value = 4383.4998
$ws.Range("L132").Value = 7123.600199999999
$ws.Range("M132").Value = -1853.4998
$ws.Range("N132").Value = -12183.6002

$ws.Range("H134").Value = 2846.3076
$ws.Range("I134").Value = 1276.4445
$ws.Range("J134").Value = 6378.5
$ws.Range("K134").Value = 3829.3335
$ws.Range("L134").Value = 19135.5
$ws.Range("M134").Value = -1294.3335
$ws.Range("N134").Value = -24205.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 2063.6956
$ws.Range("I134").Value = 1235
$ws.Range("K134").Value = 3705
$ws.Range("M134").Value = 1365

$ws.Range("H140").Value = 3711.1155
$ws.Range("I140").Value = 2638.2778
$ws.Range("J140").Value = 6125
$ws.Range("K140").Value = 7914.8334
$ws.Range("L140").Value = 18375
$ws.Range("M140").Value = -2734.8334
$ws.Range("N140").Value = -28735

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 8335460
$ws.Range("I122").Value = 9092421
$ws.Range("J122").Value = 8888
$ws.Range("K122").Value = 27277263
$ws.Range("L122").Value = 26664
$ws.Range("M122").Value = -27274813
$ws.Range("N122").Value = -31564

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1838.79
$ws.Range("I132").Value = 1840.5385
$ws.Range("J132").Value = 1821.1111
$ws.Range("K132").Value = 5521.6155
$ws.Range("L132").Value = 5463.3333
$ws.Range("M132").Value = -2991.6155
$ws.Range("N132").Value = -10523.3333

$ws.Range("H136").Value = 2000.3334
$ws.Range("J136").Value = 4215.769
$ws.Range("L136").Value = 12647.307
$ws.Range("N136").Value = -17747.307

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 559.7
$ws.Range("I107").Value = 533
$ws.Range("J107").Value = 800
$ws.Range("K107").Value = 1599
$ws.Range("L107").Value = 2400
$ws.Range("M107").Value = 321
$ws.Range("N107").Value = -6240

$ws.Range("H113").Value = 847.7857
$ws.Range("I113").Value = 854.8889
$ws.Range("J113").Value = 835
$ws.Range("K113").Value = 2564.6667
$ws.Range("L113").Value = 2505
$ws.Range("N113").Value = -6845
$ws.Range("M113").Value = -394.6667000000002

$ws.Range("H136").Value = 3771.6667
$ws.Range("I136").Value = 3860.3872
$ws.Range("J136").Value = 3575.2144
$ws.Range("K136").Value = 11581.1616
$ws.Range("L136").Value = 10725.6432
$ws.Range("M136").Value = -9031.161599999999
$ws.Range("N136").Value = -15825.6432
